# "Generate Report for Handback"
#
# The handback pass for 054acc71-143d-4577-a902-949bf83e8973.md completed:
#   * Status flips from "Ready for handoff" -> "Handed back: in sync with en-US"
#     (this shared string is used on Overview + both language sheets, so
#     updating every cell that shows it keeps them all in sync).
#   * Each language sheet gets two new columns filled in for both data rows:
#       F - Latest Target File   (the source .md file name)
#       G - Latest Handback File (the handed-back .xlf file name)
#     populated via Hyperlinks.Add so the cell carries both the display text
#     and the external link, matching how the existing A/B/D hyperlink cells
#     were produced.
#   * de-de's "Latest Handback DateTime" (H) moves off the zero-date sentinel
#     to the real handback timestamps.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdName  = "054acc71-143d-4577-a902-949bf83e8973.md"
$zhXlf   = "054acc71-143d-4577-a902-949bf83e8973.68745c589af93a2397fe89c1dff5428548af71d5.zh-cn.xlf"
$deXlf   = "054acc71-143d-4577-a902-949bf83e8973.68745c589af93a2397fe89c1dff5428548af71d5.de-de.xlf"

$mdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/76c0e4e70b5fdd1a01a3809a4e2d97c7b3740a15/e2e/$mdName"
$zhXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e25b13e783dcc48feec3afd7740bd428591f923/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/$zhXlf"
$deXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c00b32bec62eb2a59c0a7c5cf51b46497a53e3f8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/$deXlf"

# ---- Overview sheet: refresh the status column for both rows -------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

# ---- zh-cn sheet -----------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$zh.Hyperlinks.Add($zh.Range("F2"), $mdUrl, "", "", $mdName)
$zh.Hyperlinks.Add($zh.Range("G2"), $zhXlfUrl, "", "", $zhXlf)
$zh.Hyperlinks.Add($zh.Range("F3"), $mdUrl, "", "", $mdName)
$zh.Hyperlinks.Add($zh.Range("G3"), $zhXlfUrl, "", "", $zhXlf)

$zh.Range("H2").Value = "2016-03-19 20:50:40"
$zh.Range("H3").Value = "2016-03-19 20:50:40"

# ---- de-de sheet -----------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

$de.Hyperlinks.Add($de.Range("F2"), $mdUrl, "", "", $mdName)
$de.Hyperlinks.Add($de.Range("G2"), $deXlfUrl, "", "", $deXlf)
$de.Hyperlinks.Add($de.Range("F3"), $mdUrl, "", "", $mdName)
$de.Hyperlinks.Add($de.Range("G3"), $deXlfUrl, "", "", $deXlf)

$de.Range("H2").Value = "2016-03-19 20:50:45"
$de.Range("H3").Value = "2016-03-19 20:50:45"

Write-Host "Handback report generated."
